$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111780621
$ws.Range("B2").Value = 56543
$ws.Range("C2").Value = 'Ovaliderad'
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 103021
$ws.Range("F2").Value = 'Talltita'
$ws.Range("G2").Value = 'Poecile montanus'
$ws.Range("H2").Value = '(Conrad von Baldenstein, 1827)'
$ws.Range("P2").Value = 'Sasskam, Lu lm'
$ws.Range("Q2").Value = 707631.1509720345
$ws.Range("R2").Value = 7397277.54798521
$ws.Range("S2").Value = 15
$ws.Range("T2").Value = 'Norrbotten'
$ws.Range("U2").Value = 'Jokkmokk'
$ws.Range("V2").Value = 'Lule lappmark'
$ws.Range("W2").Value = 'Jokkmokk'
$ws.Range("Y2").Value = '''2023-08-29'
$ws.Range("Z2").Value = '00:00'
$ws.Range("AA2").Value = '''2023-08-29'
$ws.Range("AB2").Value = '00:00'
$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AG2").Value = $false
$ws.Range("AW2").Value = 'Amanda Tas'
$ws.Range("AX2").Value = 'Amanda Tas'

# Row 3
$ws.Range("A3").Value = 111780624
$ws.Range("B3").Value = 95532
$ws.Range("C3").Value = 'Ovaliderad'
$ws.Range("D3").Value = 'LC'
$ws.Range("E3").Value = 221945
$ws.Range("F3").Value = 'Revlummer'
$ws.Range("G3").Value = 'Lycopodium annotinum'
$ws.Range("H3").Value = 'L.'
$ws.Range("P3").Value = 'Sasskam, Lu lm'
$ws.Range("Q3").Value = 707600.9335272597
$ws.Range("R3").Value = 7397313.141869167
$ws.Range("S3").Value = 15
$ws.Range("T3").Value = 'Norrbotten'
$ws.Range("U3").Value = 'Jokkmokk'
$ws.Range("V3").Value = 'Lule lappmark'
$ws.Range("W3").Value = 'Jokkmokk'
$ws.Range("Y3").Value = '''2023-08-29'
$ws.Range("Z3").Value = '00:00'
$ws.Range("AA3").Value = '''2023-08-29'
$ws.Range("AB3").Value = '00:00'
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AW3").Value = 'Amanda Tas'
$ws.Range("AX3").Value = 'Amanda Tas'

# Row 4
$ws.Range("A4").Value = 111780628
$ws.Range("B4").Value = 78604
$ws.Range("C4").Value = 'Ovaliderad'
$ws.Range("D4").Value = 'LC'
$ws.Range("E4").Value = 6461
$ws.Range("F4").Value = 'Norrlandslav'
$ws.Range("G4").Value = 'Nephroma arcticum'
$ws.Range("H4").Value = '(L.) Torss.'
$ws.Range("P4").Value = 'Sasskam, Lu lm'
$ws.Range("Q4").Value = 707614.4806057075
$ws.Range("R4").Value = 7397255.163644295
$ws.Range("S4").Value = 15
$ws.Range("T4").Value = 'Norrbotten'
$ws.Range("U4").Value = 'Jokkmokk'
$ws.Range("V4").Value = 'Lule lappmark'
$ws.Range("W4").Value = 'Jokkmokk'
$ws.Range("Y4").Value = '''2023-08-29'
$ws.Range("Z4").Value = '00:00'
$ws.Range("AA4").Value = '''2023-08-29'
$ws.Range("AB4").Value = '00:00'
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AW4").Value = 'Amanda Tas'
$ws.Range("AX4").Value = 'Amanda Tas'

# Row 5
$ws.Range("A5").Value = 111780627
$ws.Range("B5").Value = 78604
$ws.Range("C5").Value = 'Ovaliderad'
$ws.Range("D5").Value = 'LC'
$ws.Range("E5").Value = 6461
$ws.Range("F5").Value = 'Norrlandslav'
$ws.Range("G5").Value = 'Nephroma arcticum'
$ws.Range("H5").Value = '(L.) Torss.'
$ws.Range("P5").Value = 'Sasskam, Lu lm'
$ws.Range("Q5").Value = 707647.2196405758
$ws.Range("R5").Value = 7397286.731778639
$ws.Range("S5").Value = 15
$ws.Range("T5").Value = 'Norrbotten'
$ws.Range("U5").Value = 'Jokkmokk'
$ws.Range("V5").Value = 'Lule lappmark'
$ws.Range("W5").Value = 'Jokkmokk'
$ws.Range("Y5").Value = '''2023-08-29'
$ws.Range("Z5").Value = '00:00'
$ws.Range("AA5").Value = '''2023-08-29'
$ws.Range("AB5").Value = '00:00'
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AW5").Value = 'Amanda Tas'
$ws.Range("AX5").Value = 'Amanda Tas'

# Row 6
$ws.Range("A6").Value = 111816118
$ws.Range("B6").Value = 78107
$ws.Range("C6").Value = 'Ovaliderad'
$ws.Range("D6").Value = 'NT'
$ws.Range("E6").Value = 6453
$ws.Range("F6").Value = 'Vedskivlav'
$ws.Range("G6").Value = 'Hertelidea botryosa'
$ws.Range("H6").Value = '(Fr.) Printzen & Kantvilas'
$ws.Range("P6").Value = 'Saskam, Lu lm'
$ws.Range("Q6").Value = 707670.4513803272
$ws.Range("R6").Value = 7397327.948038339
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = 'Norrbotten'
$ws.Range("U6").Value = 'Jokkmokk'
$ws.Range("V6").Value = 'Lule lappmark'
$ws.Range("W6").Value = 'Jokkmokk'
$ws.Range("Y6").Value = '''2023-08-22'
$ws.Range("Z6").Value = '00:00'
$ws.Range("AA6").Value = '''2023-08-22'
$ws.Range("AB6").Value = '00:00'
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AW6").Value = 'Cecilia Lundin'
$ws.Range("AX6").Value = 'Cecilia Lundin'

# Row 7
$ws.Range("A7").Value = 111816142
$ws.Range("B7").Value = 78604
$ws.Range("C7").Value = 'Ovaliderad'
$ws.Range("D7").Value = 'LC'
$ws.Range("E7").Value = 6461
$ws.Range("F7").Value = 'Norrlandslav'
$ws.Range("G7").Value = 'Nephroma arcticum'
$ws.Range("H7").Value = '(L.) Torss.'
$ws.Range("P7").Value = 'Saskam, Lu lm'
$ws.Range("Q7").Value = 707613.3456041727
$ws.Range("R7").Value = 7397270.22663033
$ws.Range("S7").Value = 10
$ws.Range("T7").Value = 'Norrbotten'
$ws.Range("U7").Value = 'Jokkmokk'
$ws.Range("V7").Value = 'Lule lappmark'
$ws.Range("W7").Value = 'Jokkmokk'
$ws.Range("Y7").Value = '''2023-08-29'
$ws.Range("Z7").Value = '00:00'
$ws.Range("AA7").Value = '''2023-08-29'
$ws.Range("AB7").Value = '00:00'
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AW7").Value = 'Cecilia Lundin'
$ws.Range("AX7").Value = 'Cecilia Lundin'

# Row 8
$ws.Range("A8").Value = 111816132
$ws.Range("B8").Value = 95532
$ws.Range("C8").Value = 'Ovaliderad'
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 221945
$ws.Range("F8").Value = 'Revlummer'
$ws.Range("G8").Value = 'Lycopodium annotinum'
$ws.Range("H8").Value = 'L.'
$ws.Range("P8").Value = 'Saskam, Lu lm'
$ws.Range("Q8").Value = 707589.6730983062
$ws.Range("R8").Value = 7397240.139162621
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = 'Norrbotten'
$ws.Range("U8").Value = 'Jokkmokk'
$ws.Range("V8").Value = 'Lule lappmark'
$ws.Range("W8").Value = 'Jokkmokk'
$ws.Range("Y8").Value = '''2023-08-22'
$ws.Range("Z8").Value = '00:00'
$ws.Range("AA8").Value = '''2023-08-22'
$ws.Range("AB8").Value = '00:00'
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AW8").Value = 'Cecilia Lundin'
$ws.Range("AX8").Value = 'Cecilia Lundin'

# Row 9
$ws.Range("A9").Value = 111816137
$ws.Range("B9").Value = 90658
$ws.Range("C9").Value = 'Ovaliderad'
$ws.Range("D9").Value = 'NT'
$ws.Range("E9").Value = 4361
$ws.Range("F9").Value = 'Orange taggsvamp'
$ws.Range("G9").Value = 'Hydnellum aurantiacum'
$ws.Range("H9").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("P9").Value = 'Saskam, Lu lm'
$ws.Range("Q9").Value = 707609.3988008115
$ws.Range("R9").Value = 7397264.348220735
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = 'Norrbotten'
$ws.Range("U9").Value = 'Jokkmokk'
$ws.Range("V9").Value = 'Lule lappmark'
$ws.Range("W9").Value = 'Jokkmokk'
$ws.Range("Y9").Value = '''2023-08-22'
$ws.Range("Z9").Value = '00:00'
$ws.Range("AA9").Value = '''2023-08-22'
$ws.Range("AB9").Value = '00:00'
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AW9").Value = 'Cecilia Lundin'
$ws.Range("AX9").Value = 'Cecilia Lundin'

# Row 10
$ws.Range("A10").Value = 111816145
$ws.Range("B10").Value = 77597
$ws.Range("C10").Value = 'Ovaliderad'
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 864
$ws.Range("F10").Value = 'Knottrig blåslav'
$ws.Range("G10").Value = 'Hypogymnia bitteri'
$ws.Range("H10").Value = '(Lynge) Ahti'
$ws.Range("P10").Value = 'Saskam, Lu lm'
$ws.Range("Q10").Value = 707626.9948496711
$ws.Range("R10").Value = 7397311.517900761
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = 'Norrbotten'
$ws.Range("U10").Value = 'Jokkmokk'
$ws.Range("V10").Value = 'Lule lappmark'
$ws.Range("W10").Value = 'Jokkmokk'
$ws.Range("Y10").Value = '''2023-08-22'
$ws.Range("Z10").Value = '00:00'
$ws.Range("AA10").Value = '''2023-08-22'
$ws.Range("AB10").Value = '00:00'
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AW10").Value = 'Cecilia Lundin'
$ws.Range("AX10").Value = 'Cecilia Lundin'

# Row 11
$ws.Range("A11").Value = 111816119
$ws.Range("B11").Value = 56543
$ws.Range("C11").Value = 'Ovaliderad'
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 103021
$ws.Range("F11").Value = 'Talltita'
$ws.Range("G11").Value = 'Poecile montanus'
$ws.Range("H11").Value = '(Conrad von Baldenstein, 1827)'
$ws.Range("P11").Value = 'Saskam, Lu lm'
$ws.Range("Q11").Value = 707595.5401507822
$ws.Range("R11").Value = 7397262.905378895
$ws.Range("S11").Value = 10
$ws.Range("T11").Value = 'Norrbotten'
$ws.Range("U11").Value = 'Jokkmokk'
$ws.Range("V11").Value = 'Lule lappmark'
$ws.Range("W11").Value = 'Jokkmokk'
$ws.Range("Y11").Value = '''2023-08-22'
$ws.Range("Z11").Value = '00:00'
$ws.Range("AA11").Value = '''2023-08-22'
$ws.Range("AB11").Value = '00:00'
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AG11").Value = $false
$ws.Range("AW11").Value = 'Cecilia Lundin'
$ws.Range("AX11").Value = 'Cecilia Lundin'

# Row 12
$ws.Range("A12").Value = 112202299
$ws.Range("B12").Value = 55611
$ws.Range("C12").Value = 'Ovaliderad'
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 102612
$ws.Range("F12").Value = 'Järpe'
$ws.Range("G12").Value = 'Tetrastes bonasia'
$ws.Range("H12").Value = '(Linnaeus, 1758)'
$ws.Range("P12").Value = 'Saskam, Lu lm'
$ws.Range("Q12").Value = 707645.8741767473
$ws.Range("R12").Value = 7397378.715239713
$ws.Range("S12").Value = 10
$ws.Range("T12").Value = 'Norrbotten'
$ws.Range("U12").Value = 'Jokkmokk'
$ws.Range("V12").Value = 'Lule lappmark'
$ws.Range("W12").Value = 'Jokkmokk'
$ws.Range("Y12").Value = '''2023-09-17'
$ws.Range("Z12").Value = '00:00'
$ws.Range("AA12").Value = '''2023-09-17'
$ws.Range("AB12").Value = '00:00'
$ws.Range("AD12").Value = $false
$ws.Range("AE12").Value = $false
$ws.Range("AG12").Value = $false
$ws.Range("AW12").Value = 'Cecilia Lundin'
$ws.Range("AX12").Value = 'Cecilia Lundin'
